$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Assignment 5 ---
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Assignment 5"

$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = 45184
$ws.Range("C6").Borders.Item(7).LineStyle = -4142
$ws.Range("C6").Borders.Item(10).LineStyle = -4142

$ws.Range("B6").Value = "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_5/Calc"
$ws.Range("B6").Font.Bold = $false
$ws.Range("B6").Font.Underline = $false
$ws.Range("B6").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("B6").Interior.Color = $ws.Range("C6").Interior.Color
$ws.Range("B6").Borders.Item(7).LineStyle = 1
$ws.Range("B6").VerticalAlignment = -4108

$ws.Rows("6").RowHeight = 40.5

# --- Row 7: Assignment 6 ---
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Assignment 6"

$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 45187
$ws.Range("C7").Borders.Item(7).LineStyle = -4142
$ws.Range("C7").Borders.Item(10).LineStyle = -4142

$ws.Range("B7").Value = "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_6"
$ws.Range("B7").Font.Bold = $false
$ws.Range("B7").Font.Underline = $false
$ws.Range("B7").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("B7").Interior.Color = $ws.Range("C7").Interior.Color
$ws.Range("B7").Borders.Item(7).LineStyle = 1
$ws.Range("B7").VerticalAlignment = -4108

$ws.Rows("7").RowHeight = 37.5

# --- Row 8: Assignment 7 ---
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Assignment 7"

$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = 45188
$ws.Range("C8").Borders.Item(7).LineStyle = -4142
$ws.Range("C8").Borders.Item(10).LineStyle = -4142

$ws.Range("B8").Value = "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_7"
$ws.Range("B8").Font.Bold = $false
$ws.Range("B8").Font.Underline = $false
$ws.Range("B8").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("B8").Interior.Color = $ws.Range("C8").Interior.Color
$ws.Range("B8").Borders.Item(7).LineStyle = 1
$ws.Range("B8").VerticalAlignment = -4108

$ws.Rows("8").RowHeight = 36.75

$ws.Range("B8").Select()
